$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1835
$ws.Range("G2").Value = 1835
$ws.Range("D3").Value = 16
$ws.Range("G3").Value = 16
